# Project's bugs.xlsx - apply commit:
# "Fix sort by Status (M + C) Fix sort by Scale (M) Fix bugs"
#
# - Mark the two "Sort by Status not working" bug rows (Merchant + Customer)
#   and the "Show Gender" bug row as Fixed.
# - Add two new bug rows (Backend/Customer + Backend/Merchant), both Not fix.
# - Update the visible scroll position / active selection on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Flip Status column (F) from "Not fix" to "Fixed" for rows 2, 5 and 30 ---
$ws.Range("F2").Value  = "Fixed"
$ws.Range("F5").Value  = "Fixed"
$ws.Range("F30").Value = "Fixed"

# --- New row 31: Backend / Customer / "can't update to table User " ---
$ws.Range("A31").Value = 30
$ws.Range("B31").Value = "Backend"
$ws.Range("C31").Value = "Customer"
$ws.Range("D31").Value = "can't update to table User "
$ws.Range("E31").Value = "Bug"
$ws.Range("F31").Value = "Not fix"
$ws.Range("G31").Value = "Cá"

# --- New row 32: Backend / Merchant / "change MerchantCompanyName?" ---
$ws.Range("A32").Value = 31
$ws.Range("B32").Value = "Backend"
$ws.Range("C32").Value = "Merchant"
$ws.Range("D32").Value = "change MerchantCompanyName?"
$ws.Range("E32").Value = "Bug?"
$ws.Range("F32").Value = "Not fix"
$ws.Range("G32").Value = "Cá"

# --- Update view state: scrolled down a couple more rows, new active cell ---
$ws.Range("H32").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1
